# Apply updated hyperparameter search run-data values to row 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.1816669225692749
$ws.Range("B2").Value = 542.2857142857143
$ws.Range("C2").Value = 5.978195488721804
$ws.Range("H2").Value = 150
$ws.Range("I2").Value = 15
$ws.Range("M2").Value = 0.05
$ws.Range("V2").Value = 5
